$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Remove paragraphs that are dropped entirely by the edit.
#    (delete from the bottom up so indices of earlier paragraphs
#    stay stable while we work)
# ---------------------------------------------------------------

# "- Other data are linked to the group within the heating system. ..."
$d.Paragraphs.Item(10).Range.Delete()

# "- Data on thermal properties marked in red are the properties ..."
$d.Paragraphs.Item(9).Range.Delete()

# "The following fields are all boolean values that tell whether ..."
$d.Paragraphs.Item(7).Range.Delete()

# ---------------------------------------------------------------
# 2. "- Flags. These fields serve ..." -> "Flags are boolean values
#    that serve ..." and "... we will implement that." -> "... we
#    will discuss that."
# ---------------------------------------------------------------
$r = $d.Paragraphs.Item(7).Range
$r.Find.Execute("- Flags. These fields serve", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Flags are boolean values that serve", 2)

$r = $d.Paragraphs.Item(7).Range
$r.Find.Execute("open an issue and we will implement that.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "open an issue and we will discuss that.", 2)

# ---------------------------------------------------------------
# 3. "- The info.json file ..." paragraph edits
# ---------------------------------------------------------------
$r = $d.Paragraphs.Item(9).Range
$r.Find.Execute("Temperature ranges must be recorded for density, specific heat capacity, thermal conductivity, adiabatic temperature change (if the material is caloric), emissivity, and other relevant properties (e.g. seebeck coefficient for thermoelectric materials, etc.). The fields", `
                 $true, $false, $false, $false, $false, $true, 1, $false, `
                 "Temperature ranges must be given for density, specific heat capacity, thermal conductivity, and emissivity. The fields", 2)

$r = $d.Paragraphs.Item(9).Range
$r.Find.Execute("and for pressure and stress, they are in bars. (See any material for example).", `
                 $true, $false, $false, $false, $false, $true, 1, $false, `
                 "and for pressure and stress, they are in kbars. For any unknown values, input “”.", 2)

# ---------------------------------------------------------------
# 4. "- Files rho.txt, cp.txt and k.txt ..." paragraph edit
# ---------------------------------------------------------------
$r = $d.Paragraphs.Item(10).Range
$r.Find.Execute("one column of 20000 values from 0 to 2000 K", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "one column of 20000 values from 0.1 to 2000 K", 2)

# ---------------------------------------------------------------
# 5. "- One or more of the above three files ..." paragraph edit
# ---------------------------------------------------------------
$r = $d.Paragraphs.Item(11).Range
$r.Find.Execute("electric in MVm (which means MV/m)), e.g.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "electric in MVm (which stands for MV/m), stress and pressure in kbar), e.g.", 2)

# ---------------------------------------------------------------
# 6. "- There are also other options ..." paragraph edit (append
#    a new trailing sentence)
# ---------------------------------------------------------------
$r = $d.Paragraphs.Item(12).Range
$r.Find.Execute("cp_0.0T_cooling.txt, cp_0.0T_heating.txt, cp_1.0T_cooling.txt, etc.", `
                 $true, $false, $false, $false, $false, $true, 1, $false, `
                 "cp_0.0T_cooling.txt, cp_0.0T_heating.txt, cp_1.0T_cooling.txt, etc. Or, we could have cp data for different values of external field without hysteresis (cpFields flag is true). Please note that the hysteresis flags tell the program which files to look for and are not necessary consistent with the actual hysteresis of the material.", 2)

# ---------------------------------------------------------------
# 7. Brand new closing paragraph about adiabatic temperature change
# ---------------------------------------------------------------
$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = "- Caloric materials will have adiabatic temperature change, which must be zero outside the range, where it is defined. It must be given for the case of application of external field and removal of external field. If necessary, give information on the range of temperatures where the caloric effect is reversible in a file in data folder."
